$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "https://codeforces.com/problemset/problem/1881/A"
$ws.Range("C9").Value = "how tpo find subset of a strig"
$ws.Range("D9").Value = "temp.found(s) != string::npos (use to checl if s is a subset of temp or not and string::npose means tring not found"
